$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# Change 1: there is a run of 9 consecutive empty "Normal" paragraphs right
# before the bold "Matrículas vehiculares" heading. Trim that down to a
# single empty paragraph (remove 8 of the 9).
# --------------------------------------------------------------------------
$headingIndex = -1
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Matrículas vehiculares*") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -gt 0) {
    # Walk backwards from the paragraph right before the heading, counting
    # consecutive empty "Normal" style paragraphs.
    $emptyCount = 0
    $j = $headingIndex - 1
    while ($j -ge 1) {
        $p = $d.Paragraphs.Item($j)
        if ($p.Range.Text -eq [char]13 -and $p.Style.NameLocal -eq "Normal") {
            $emptyCount = $emptyCount + 1
            $j = $j - 1
        } else {
            break
        }
    }

    # Keep exactly one empty paragraph; delete the rest (the ones closest to
    # the start of that run, i.e. the earliest indices).
    if ($emptyCount -gt 1) {
        $toDelete = $emptyCount - 1
        $firstEmptyIndex = $headingIndex - $emptyCount
        $lastDeleteIndex = $firstEmptyIndex + $toDelete - 1

        for ($k = $lastDeleteIndex; $k -ge $firstEmptyIndex; $k--) {
            $pd = $d.Paragraphs.Item($k)
            $pd.Range.Delete()
        }
    }
}

# --------------------------------------------------------------------------
# Change 2: after the paragraph that reads "AA-00-000", add two new pattern
# lines ("0-A-00-AA" and "0-A-0-AAA") followed by one blank paragraph before
# the existing "Estos patrones..." paragraph.
# --------------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*AA-00-000*") {
        $p.Range.InsertParagraphAfter()
        $newPara1 = $d.Paragraphs.Item($i + 1)
        $newPara1.Range.Text = "0-A-00-AA"

        $newPara1.Range.InsertParagraphAfter()
        $newPara2 = $d.Paragraphs.Item($i + 2)
        $newPara2.Range.Text = "0-A-0-AAA"

        $newPara2.Range.InsertParagraphAfter()
        break
    }
}
